$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets 1-4 ("Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)") share the same
# layout: a header row (B1:E1 = years) and a "Fonte/Tecnologia" label column
# (A2:A12). Add the missing header label in A1 (reusing the existing header
# style from B1), accent-fix a few of the technology names, and drop the
# heavy "bold/border" style from the label cells now that the real header
# row exists.
# ---------------------------------------------------------------------------
$dataSheetNames = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($name in $dataSheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # New header label for column A - give it the same look as the other
    # header cells (B1:E1) by copying their formatting over.
    $ws.Range("A1").Value = "Fonte/Tecnologia"
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    # Corrected (accented) technology names.
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."

    # Remove the bold/border style from the label column (rows 2-12) now
    # that it is no longer needed there.
    for ($r = 2; $r -le 12; $r++) {
        $ws.Cells.Item($r, 1).Style = "Normal"
    }
}

# ---------------------------------------------------------------------------
# Sheet 5 ("Emissoes Totais (MtCO2eq)"): add a "Período" header label, fix
# the accented row labels, drop their style, and remove the now-unused
# "Teto" row.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")

$ws5.Range("A1").Value = "Período"
$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A2").Style = "Normal"

$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A3").Style = "Normal"

# Drop row 4 ("Teto") entirely.
$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# Sheet 6 ("Custo Total (bilhões de R$)"): add the "Tipo Expansão" header
# label, rename the year header, fix the accented row labels, drop their
# style, and update the cost values.
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")

$ws6.Range("A1").Value = "Tipo Expansão"
$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# B1 becomes the text label "2015" (not a number) - force text entry, then
# restore the original "General" number formatting (taken from A1, which now
# carries the same header style) without disturbing the stored text value.
$ws6.Range("B1").NumberFormat = "@"
$ws6.Range("B1").Value = "2015"
$ws6.Range("A1").Copy()
$ws6.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("A2").Style = "Normal"
$ws6.Range("B2").Value = 475

$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("A3").Style = "Normal"
$ws6.Range("B3").Value = 99
